# Actualización desde MV -datos-
# Update values in Sheet1 to reflect refreshed data from the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C26").Value = 80.90000000000001

$ws.Range("D34").Value = 1.3

$ws.Range("D46").Value = 0.2

$ws.Range("C50").Value = 103

$ws.Range("C56").Value = 105.7
$ws.Range("D56").Value = 0.3

$ws.Range("C57").Value = 106

$ws.Range("D59").Value = 1.1

$ws.Range("C60").Value = 108.1
$ws.Range("D60").Value = 1.7

$ws.Range("C63").Value = 111.4

$ws.Range("C64").Value = 110.9
$ws.Range("D64").Value = -0.4

$ws.Range("D65").Value = 1.2

$ws.Range("C66").Value = 112
$ws.Range("D66").Value = -0.2

$ws.Range("C67").Value = 113.3
$ws.Range("D67").Value = 1.2

$ws.Range("C68").Value = 113.9
$ws.Range("D68").Value = 0.6

$ws.Range("C69").Value = 109.6

$ws.Range("C71").Value = 97.8

$ws.Range("C72").Value = 102.8
$ws.Range("D72").Value = 5.1

$ws.Range("C73").Value = 109.5
$ws.Range("D73").Value = 6.6

$ws.Range("B74").Value = 110.5
$ws.Range("C74").Value = 113.3
$ws.Range("D74").Value = 3.5

$ws.Range("B75").Value = 115.1
$ws.Range("C75").Value = 114.6
$ws.Range("D75").Value = 1.2

$wb.Save()
